$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Enhancements 56 & 57": split the single "Fecha inicio / Fecha fin" column
# into two dedicated columns by inserting a new blank column before E.
# Everything that used to live in E..I (Fecha inicio, Fecha fin, Contrato
# renta, Contrato servicios, Cliente + the per-column "date"/"datetime"
# hint cells in row 3) shifts one column to the right, into F..J, and the
# new column D:E date-range block gets a matching wider width.
$ws.Columns("E:E").Insert()

# Give the now-paired D:E columns (the old "Fecha" block) a consistent
# width like the rest of the wide text columns.
$ws.Columns("D:E").ColumnWidth = 47.6

# Restore/extend the AutoFilter over the new last column (J) - inserting a
# column does not automatically grow the filter range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A2:J2").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# AutoFilter range above.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Contratos!_FilterDatabase") {
        $n.RefersTo = "=Contratos!`$A`$2:`$J`$2"
    }
}
